# Adds the "iNaturalist ID" field to the voucher template header table
# (Sheet1, rows 1-2), and removes the old QR-image-code fields
# (site_image_code / specimen_image_code) that it replaces.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before AW (ProcessingStatus and everything after it
#    shifts one column to the right) and give the new header its label.
$ws.Columns("AW:AW").Insert()
$ws.Range("AW1").Value = "iNaturalist ID"
$ws.Range("AW2").Value = ""

# 2) The two trailing "QR image method" columns (site_image_code /
#    specimen_image_code), which used to be BC:BD, are now BD:BE after the
#    insert above - remove them entirely.
$ws.Columns("BD:BE").Delete()

# 3) Resize the new column to fit its header text, like Excel would when a
#    column is added and used.
$ws.Columns("AW:AW").AutoFit()

# 4) Match the resulting selection/scroll position left behind by the edit.
$ws.Range("AS1").Select()
$ws.Range("AW2").Select()
